# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the newer scrape referenced in the commit message.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, new value)
$updates = @{
    "展览" = @{
        "F5"  = 15887
        "F9"  = 15488
        "F11" = 9100
        "F15" = 110
        "F18" = 209
        "F21" = 570
        "F24" = 64
        "F29" = 478
        "F34" = 51
        "F35" = 263
        "F39" = 5594
    }
    "全部类型" = @{
        "F5"  = 15887
        "F9"  = 15488
        "F11" = 9100
        "F15" = 110
        "F18" = 209
        "F21" = 570
        "F24" = 64
        "F29" = 478
        "F36" = 51
        "F37" = 263
        "F41" = 5594
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}

$wb.Save()
